$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112370039
$ws.Range("B2").Value = 77650
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = 'Garnlav'
$ws.Range("G2").Value = 'Alectoria sarmentosa'
$ws.Range("H2").Value = '(Ach.) Ach.'
$ws.Range("Q2").Value = 469889
$ws.Range("R2").Value = 7039252
$ws.Range("AC2").Value = ""

# Row 3
$ws.Range("A3").Value = 112370010
$ws.Range("Q3").Value = 469821
$ws.Range("R3").Value = 7039232
$ws.Range("AC3").Value = 'ringhack äldre'

# Row 4
$ws.Range("A4").Value = 112370025
$ws.Range("B4").Value = 56446
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = 'Spillkråka'
$ws.Range("G4").Value = 'Dryocopus martius'
$ws.Range("H4").Value = '(Linnaeus, 1758)'
$ws.Range("Q4").Value = 469994
$ws.Range("R4").Value = 7039247
$ws.Range("AC4").Value = 'hack'

# Row 5
$ws.Range("A5").Value = 112370038
$ws.Range("B5").Value = 77650
$ws.Range("Q5").Value = 469826
$ws.Range("R5").Value = 7039235

# Row 6
$ws.Range("A6").Value = 112370008
$ws.Range("Q6").Value = 469854
$ws.Range("R6").Value = 7039173
$ws.Range("AC6").Value = 'ringhack äldre'

# Row 7
$ws.Range("A7").Value = 112370009
$ws.Range("B7").Value = 56430
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 100109
$ws.Range("F7").Value = 'Tretåig hackspett'
$ws.Range("G7").Value = 'Picoides tridactylus'
$ws.Range("H7").Value = '(Linnaeus, 1758)'
$ws.Range("Q7").Value = 469795
$ws.Range("R7").Value = 7039224
$ws.Range("AC7").Value = 'ringhack'

# Row 8
$ws.Range("A8").Value = 112370041
$ws.Range("B8").Value = 90835
$ws.Range("D8").Value = 'LC'
$ws.Range("E8").Value = 5964
$ws.Range("F8").Value = 'Fjällig taggsvamp s.str.'
$ws.Range("G8").Value = 'Sarcodon imbricatus s.str.'
$ws.Range("H8").Value = '(L.:Fr.) P.Karst.'
$ws.Range("Q8").Value = 469969
$ws.Range("R8").Value = 7039249
$ws.Range("AC8").Value = ""

# Row 9
$ws.Range("A9").Value = 112370007
$ws.Range("B9").Value = 56430
$ws.Range("E9").Value = 100109
$ws.Range("F9").Value = 'Tretåig hackspett'
$ws.Range("G9").Value = 'Picoides tridactylus'
$ws.Range("Q9").Value = 469863
$ws.Range("R9").Value = 7039172
$ws.Range("AC9").Value = 'ringhack'

# Row 10
$ws.Range("B10").Value = 83086

# Row 11
$ws.Range("B11").Value = 85850

# Row 12
$ws.Range("A12").Value = 112370045
$ws.Range("B12").Value = 90799
$ws.Range("E12").Value = 1968
$ws.Range("F12").Value = 'Grantaggsvamp'
$ws.Range("G12").Value = 'Bankera violascens'
$ws.Range("H12").Value = '(Alb. & Schwein. : Fr.) Pouzar'
$ws.Range("Q12").Value = 470217
$ws.Range("R12").Value = 7038987
$ws.Range("Y12").NumberFormat = "@"
$ws.Range("Y12").Value = '2023-09-28'
$ws.Range("Y12").ClearFormats()
$ws.Range("AA12").NumberFormat = "@"
$ws.Range("AA12").Value = '2023-09-28'
$ws.Range("AA12").ClearFormats()

# Row 13
$ws.Range("A13").Value = 112370040
$ws.Range("B13").Value = 77650
$ws.Range("E13").Value = 6425
$ws.Range("F13").Value = 'Garnlav'
$ws.Range("G13").Value = 'Alectoria sarmentosa'
$ws.Range("H13").Value = '(Ach.) Ach.'
$ws.Range("Q13").Value = 470143
$ws.Range("R13").Value = 7039263

# Row 14
$ws.Range("A14").Value = 112370037
$ws.Range("B14").Value = 89571
$ws.Range("E14").Value = 5432
$ws.Range("F14").Value = 'Granticka'
$ws.Range("G14").Value = 'Porodaedalea chrysoloma'
$ws.Range("H14").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q14").Value = 470039
$ws.Range("R14").Value = 7039048
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = '2023-09-26'
$ws.Range("Y14").ClearFormats()
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = '2023-09-26'
$ws.Range("AA14").ClearFormats()
